$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 35, shifting existing rows 35:84 down to 36:85.
$ws.Rows("35:35").Insert()

# Populate the newly inserted row 35 with the new data record.
$ws.Cells.Item(35, 1).Value = 5
$ws.Cells.Item(35, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(35, 3).Value = "Maule"
$ws.Cells.Item(35, 4).Value = 44771
$ws.Cells.Item(35, 5).Value = 7
$ws.Cells.Item(35, 6).Value = 100112013
$ws.Cells.Item(35, 7).Value = "Alcachofa"
$ws.Cells.Item(35, 8).Value = "Madrigal"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 400
$ws.Cells.Item(35, 11).Value = 13000
$ws.Cells.Item(35, 12).Value = 13000
$ws.Cells.Item(35, 13).Value = 13000
$ws.Cells.Item(35, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(35, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(35, 16).Value = 325
$ws.Cells.Item(35, 17).Value = 40
$ws.Cells.Item(35, 18).Value = "Hortaliza"
